$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.524.62"
$ws.Range("E2").Value = "  +5.70%  "
$ws.Range("D3").Value = "1.723.59"
$ws.Range("E3").Value = "  +4.66%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'225.87"
$ws.Range("E5").Value = "  +3.71%  "
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.2675"
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("D9").Value = "'0.06590"
$ws.Range("E9").Value = "  +4.81%  "
$ws.Range("D10").Value = "'21.65"
$ws.Range("E10").Value = "  +6.80%  "
$ws.Range("D11").Value = "'0.07712"
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "1.722.77"
$ws.Range("E13").Value = "  +4.48%  "
$ws.Range("D14").Value = "1.960.98"
$ws.Range("E14").Value = "  +4.66%  "
$ws.Range("D15").Value = "'0.5833"
$ws.Range("E15").Value = "  +5.09%  "
$ws.Range("D16").Value = "0.0₅8283"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "'68.01"
$ws.Range("E17").Value = "  +4.74%  "
$ws.Range("D18").Value = "27.523.75"
$ws.Range("E18").Value = "  +5.77%  "
$ws.Range("D19").Value = "'219.47"
$ws.Range("E19").Value = "  +15.25%  "
$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'4.738"
$ws.Range("E21").Value = "  +3.28%  "
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").Value = "'6.081"
$ws.Range("E23").Value = "  +3.49%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'145.98"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "'1.728"
$ws.Range("E26").Value = "  +14.42%  "
$ws.Range("D27").Value = "'0.1235"
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").Value = "'7.405"
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("D29").Value = "'16.57"
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("D30").Value = "'0.05541"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("D31").Value = "'1.305"
$ws.Range("E31").Value = "  +3.23%  "
$ws.Range("D32").Value = "'3.566"
$ws.Range("E32").Value = "  +3.81%  "
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("E34").Value = "  +8.20%  "
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("D36").Value = "'0.9669"
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("D37").Value = "'2.424"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("E38").Value = "  +7.37%  "
$ws.Range("E39").Value = "  +5.40%  "
$ws.Range("D40").Value = "'5.913"
$ws.Range("E40").Value = "  +2.55%  "
$ws.Range("D41").Value = "'0.8579"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("D42").Value = "1.055.72"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'101.32"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").Value = "1.867.21"
$ws.Range("E45").Value = "  +4.54%  "
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("D47").Value = "'59.00"
$ws.Range("E47").Value = "  +3.74%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.4454"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.204"
$ws.Range("E49").Value = "  +4.66%  "
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'0.05234"
$ws.Range("E51").Value = "  +2.28%  "
